$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "1.000")
# are preserved verbatim instead of being parsed into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.683.80'
$ws.Range("E2").Value = '  +0.72%  '
$ws.Range("D3").Value = '1.923.37'
$ws.Range("E3").Value = '  -0.65%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.67%  '
$ws.Range("D5").Value = '335.34'
$ws.Range("E5").Value = '  -1.27%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("D7").Value = '0.4671'
$ws.Range("E7").Value = '  -1.41%  '
$ws.Range("D8").Value = '0.4148'
$ws.Range("E8").Value = '  +1.19%  '
$ws.Range("D9").Value = '48.18'
$ws.Range("E9").Value = '  +0.77%  '
$ws.Range("D10").Value = '0.08069'
$ws.Range("E10").Value = '  -1.55%  '
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").Value = '22.32'
$ws.Range("E12").Value = '  -0.72%  '
$ws.Range("D13").Value = '1.936.35'
$ws.Range("E13").Value = '  -0.55%  '
$ws.Range("D14").Value = '6.013'
$ws.Range("E14").Value = '  -1.50%  '
$ws.Range("D15").Value = '7.189'
$ws.Range("E15").Value = '  -1.76%  '
$ws.Range("D16").Value = '89.84'
$ws.Range("E16").Value = '  -1.65%  '
$ws.Range("E17").Value = '  -0.73%  '
$ws.Range("D18").Value = '0.00001039'
$ws.Range("E18").Value = '  -1.20%  '
$ws.Range("D19").Value = '0.06609'
$ws.Range("E19").Value = '  -1.06%  '
$ws.Range("D20").Value = '17.82'
$ws.Range("E20").Value = '  -0.29%  '
$ws.Range("D21").Value = '0.9984'
$ws.Range("E21").Value = '  -0.40%  '
$ws.Range("D22").Value = '29.638.56'
$ws.Range("E22").Value = '  +0.51%  '
$ws.Range("D23").Value = '5.555'
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("E24").Value = '  +3.92%  '
$ws.Range("E25").Value = '  -3.38%  '
$ws.Range("D26").Value = '2.140.40'
$ws.Range("E26").Value = '  -1.39%  '
$ws.Range("D27").Value = '157.52'
$ws.Range("E27").Value = '  -2.15%  '
$ws.Range("D28").Value = '19.95'
$ws.Range("E28").Value = '  -0.61%  '
$ws.Range("D29").Value = '2.154'
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("D30").Value = '5.738'
$ws.Range("E30").Value = '  +2.59%  '
$ws.Range("D31").Value = '118.00'
$ws.Range("E31").Value = '  -3.49%  '
$ws.Range("D32").Value = '1.049'
$ws.Range("E32").Value = '  +4.86%  '
$ws.Range("D33").Value = '0.09459'
$ws.Range("E33").Value = '  -1.56%  '
$ws.Range("D34").Value = '1.434'
$ws.Range("E34").Value = '  -1.43%  '
$ws.Range("D35").Value = '5.436'
$ws.Range("E35").Value = '  -0.38%  '
$ws.Range("E36").Value = '  -3.54%  '
$ws.Range("D37").Value = '0.06149'
$ws.Range("E37").Value = '  -1.17%  '
$ws.Range("D38").Value = '0.02275'
$ws.Range("E38").Value = '  -1.14%  '
$ws.Range("D39").Value = '8.472'
$ws.Range("E39").Value = '  +0.55%  '
$ws.Range("D40").Value = '1.176'
$ws.Range("E40").Value = '  -0.70%  '
$ws.Range("D41").Value = '0.5921'
$ws.Range("E41").Value = '  -1.87%  '
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").Value = '0.1846'
$ws.Range("E43").Value = '  -2.10%  '
$ws.Range("D44").Value = '10.27'
$ws.Range("E44").Value = '  -3.26%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").Value = '2.342'
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("B46").Value = 'WEMIXTOKEN'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '1.242'
$ws.Range("E46").Value = '  -2.39%  '
$ws.Range("D47").Value = '0.07540'
$ws.Range("E47").Value = '  +2.01%  '
$ws.Range("D48").Value = '0.5595'
$ws.Range("E48").Value = '  -1.45%  '
$ws.Range("D49").Value = '12.14'
$ws.Range("E49").Value = '  -2.90%  '
$ws.Range("D50").Value = '1.946'
$ws.Range("E50").Value = '  -1.45%  '
$ws.Range("D51").Value = '113.19'
$ws.Range("E51").Value = '  +0.83%  '

# Restore default styling on column D (clears the temporary text format,
# keeping the cell style index unchanged at 0 like the rest of the sheet).
$ws.Range("D2:D51").Style = "Normal"
